# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets,
# matching the data refresh recorded in the commit "Update gh-pages to
# output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Row (key) -> new F-column value, per worksheet.
$updatesByPlanSheet = @{
    2  = 15003
    3  = 18970
    14 = 151
    15 = 217
    16 = 60
    17 = 1456
    20 = 97
    22 = 7902
    25 = 2
    26 = 62
    27 = 1241
    29 = 6042
    31 = 71
    32 = 168
    34 = 279
    35 = 5410
    36 = 124
    37 = 8
}

$updatesByAllTypesSheet = @{
    2  = 15003
    3  = 18970
    14 = 151
    15 = 217
    16 = 60
    17 = 1456
    21 = 97
    23 = 7902
    26 = 2
    27 = 62
    28 = 1241
    32 = 6042
    34 = 71
    35 = 168
    37 = 279
    38 = 5410
    39 = 124
    40 = 8
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updatesByPlanSheet.Keys) {
    $ws1.Range("F$row").Value = $updatesByPlanSheet[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updatesByAllTypesSheet.Keys) {
    $ws4.Range("F$row").Value = $updatesByAllTypesSheet[$row]
}
